$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "Leaves - random orientation, change colour, fall to ground" feature
# is now fully implemented: mark goes from 5 to full marks (10) and is
# flagged DONE in the adjacent notes column.
$ws.Range("D11").Value = 10
$ws.Range("E11").Value = "DONE"
$ws.Range("E11").Font.Color = $ws.Range("E4").Font.Color

# Match author's final cursor position recorded in the sheet view.
$ws.Range("F15").Select()

$wb.Save()
